$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.700.11"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3
$ws.Range("D3").Value = "3.445.83"
$ws.Range("E3").Value = "  +2.17%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.30%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "3.445.98"
$ws.Range("E8").Value = "  +2.19%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "

# Row 10
$ws.Range("E10").Value = "  +3.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.393"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("D13").Value = "4.029.42"
$ws.Range("E13").Value = "  +2.08%  "

# Row 14
$ws.Range("E14").Value = "  -1.02%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "

# Row 17
$ws.Range("D17").Value = "3.445.48"
$ws.Range("E17").Value = "  +2.19%  "

# Row 18
$ws.Range("D18").Value = "61.718.54"
$ws.Range("E18").Value = "  +0.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.52%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.18%  "

# Row 23
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.561"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.29%  "

# Row 24
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.565.89"
$ws.Range("E24").Value = "  +1.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.27%  "

# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.72%  "

# Row 28
$ws.Range("E28").Value = "  +10.61%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.36%  "

# Row 30
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.02%  "

# Row 34
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.33%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "166.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0795"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.21%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.791"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.34%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.88%  "

# Row 44
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "

# Row 48
$ws.Range("D48").Value = "2.658.69"
$ws.Range("E48").Value = "  +13.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.36%  "
